$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The latitude/longitude for id=10001 (rows 2-4) should all match row 2's values.
# Row 2 already holds the correct lat/long; copy it down into rows 3 and 4 so the
# duplicate text values (" -6.453276" / " -6.453277") are no longer used anywhere
# and can be dropped from the shared-string table.
$ws.Range("G2").Copy($ws.Range("G3"))
$ws.Range("H2").Copy($ws.Range("H3"))
$ws.Range("G2").Copy($ws.Range("G4"))
$ws.Range("H2").Copy($ws.Range("H4"))

# number_of_kiosks (column L): the first location (id=10001, rows 2-4) now has 3
# kiosks, and every other location (rows 5-46) has 2 kiosks.
$ws.Range("L2").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("L4").Value = 3

for ($r = 5; $r -le 46; $r++) {
    $ws.Cells.Item($r, 12).Value = 2
}

$ws.Range("L28").Select() | Out-Null
